$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: cardholder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3: card number (keep as text so it is not reformatted as a number) and last name
$fmtB3 = $ws.Range("B3").NumberFormat
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B3").NumberFormat = $fmtB3
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance line (date in heading text)
$ws.Range("D5").Value = "KONTOSTAND AM 01.01.2025"

# Row 6: transaction 1
$ws.Range("B6").Value = "03.01."
$ws.Range("C6").Value = "04.01."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "25,25-"

# Row 7: transaction 2
$ws.Range("B7").Value = "05.01."
$ws.Range("C7").Value = "06.01."
$ws.Range("D7").Value = "PAYPAL ERGVGD"
$ws.Range("E7").Value = "77,95-"

# Row 8: transaction 3
$ws.Range("B8").Value = "07.01."
$ws.Range("C8").Value = "08.01."
$ws.Range("D8").Value = "KARTENZ./07.01 ALDI SUED RO"
$ws.Range("E8").Value = "52,14-"

# Row 12: closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 10.01.2025"
$ws.Range("E12").Value = "155,34-"

# Row 13: next billing date note
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 20.01.2025"
